$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds dates formatted as dd/mm/yyyy TEXT (not real Excel dates).
# Set the number format to Text up front so the assignment below doesn't get
# auto-converted into date serial numbers by Excel.
$ws.Range("B2:B31").NumberFormat = "@"

$arr = New-Object 'object[,]' 31,6
$arr[0,0] = 'Região'
$arr[0,1] = 'Ano'
$arr[0,2] = 'Variável'
$arr[0,3] = 'Valor'
$arr[0,4] = 'Posição relativamente às demais UF'
$arr[0,5] = 'Faltam dados para todos os Estados'
$arr[1,0] = 'Brasil'
$arr[1,1] = '01/01/2015'
$arr[1,2] = 'Homicídio doloso'
$arr[1,3] = 27.28839008256385
$arr[1,4] = $null
$arr[1,5] = $true
$arr[2,0] = 'Brasil'
$arr[2,1] = '01/01/2016'
$arr[2,2] = 'Homicídio doloso'
$arr[2,3] = 30.05887666703984
$arr[2,4] = $null
$arr[2,5] = $true
$arr[3,0] = 'Brasil'
$arr[3,1] = '01/01/2017'
$arr[3,2] = 'Homicídio doloso'
$arr[3,3] = 32.89455525672006
$arr[3,4] = $null
$arr[3,5] = $false
$arr[4,0] = 'Brasil'
$arr[4,1] = '01/01/2018'
$arr[4,2] = 'Homicídio doloso'
$arr[4,3] = 29.92000675669447
$arr[4,4] = $null
$arr[4,5] = $false
$arr[5,0] = 'Brasil'
$arr[5,1] = '01/01/2019'
$arr[5,2] = 'Homicídio doloso'
$arr[5,3] = 22.78123160802507
$arr[5,4] = $null
$arr[5,5] = $false
$arr[6,0] = 'Brasil'
$arr[6,1] = '01/01/2020'
$arr[6,2] = 'Homicídio doloso'
$arr[6,3] = 23.95462631156656
$arr[6,4] = $null
$arr[6,5] = $false
$arr[7,0] = 'Brasil'
$arr[7,1] = '01/01/2021'
$arr[7,2] = 'Homicídio doloso'
$arr[7,3] = 22.45847761979459
$arr[7,4] = $null
$arr[7,5] = $false
$arr[8,0] = 'Brasil'
$arr[8,1] = '01/01/2022'
$arr[8,2] = 'Homicídio doloso'
$arr[8,3] = 21.82561212792709
$arr[8,4] = $null
$arr[8,5] = $false
$arr[9,0] = 'Brasil'
$arr[9,1] = '01/01/2023'
$arr[9,2] = 'Homicídio doloso'
$arr[9,3] = 20.91732437986478
$arr[9,4] = $null
$arr[9,5] = $false
$arr[10,0] = 'Brasil'
$arr[10,1] = '01/01/2024'
$arr[10,2] = 'Homicídio doloso'
$arr[10,3] = 18.61443649398527
$arr[10,4] = $null
$arr[10,5] = $false
$arr[11,0] = 'Nordeste'
$arr[11,1] = '01/01/2015'
$arr[11,2] = 'Homicídio doloso'
$arr[11,3] = 38.83151646101115
$arr[11,4] = $null
$arr[11,5] = $true
$arr[12,0] = 'Nordeste'
$arr[12,1] = '01/01/2016'
$arr[12,2] = 'Homicídio doloso'
$arr[12,3] = 40.38565884059511
$arr[12,4] = $null
$arr[12,5] = $true
$arr[13,0] = 'Nordeste'
$arr[13,1] = '01/01/2017'
$arr[13,2] = 'Homicídio doloso'
$arr[13,3] = 41.76396000036527
$arr[13,4] = $null
$arr[13,5] = $false
$arr[14,0] = 'Nordeste'
$arr[14,1] = '01/01/2018'
$arr[14,2] = 'Homicídio doloso'
$arr[14,3] = 35.03400939047036
$arr[14,4] = $null
$arr[14,5] = $false
$arr[15,0] = 'Nordeste'
$arr[15,1] = '01/01/2019'
$arr[15,2] = 'Homicídio doloso'
$arr[15,3] = 26.49153069038045
$arr[15,4] = $null
$arr[15,5] = $false
$arr[16,0] = 'Nordeste'
$arr[16,1] = '01/01/2020'
$arr[16,2] = 'Homicídio doloso'
$arr[16,3] = 31.7167489266386
$arr[16,4] = $null
$arr[16,5] = $false
$arr[17,0] = 'Nordeste'
$arr[17,1] = '01/01/2021'
$arr[17,2] = 'Homicídio doloso'
$arr[17,3] = 28.49517714923115
$arr[17,4] = $null
$arr[17,5] = $false
$arr[18,0] = 'Nordeste'
$arr[18,1] = '01/01/2022'
$arr[18,2] = 'Homicídio doloso'
$arr[18,3] = 27.52826693045831
$arr[18,4] = $null
$arr[18,5] = $false
$arr[19,0] = 'Nordeste'
$arr[19,1] = '01/01/2023'
$arr[19,2] = 'Homicídio doloso'
$arr[19,3] = 26.11187036839209
$arr[19,4] = $null
$arr[19,5] = $false
$arr[20,0] = 'Nordeste'
$arr[20,1] = '01/01/2024'
$arr[20,2] = 'Homicídio doloso'
$arr[20,3] = 24.76353685797129
$arr[20,4] = $null
$arr[20,5] = $false
$arr[21,0] = 'Sergipe'
$arr[21,1] = '01/01/2015'
$arr[21,2] = 'Homicídio doloso'
$arr[21,3] = 53.95512251106057
$arr[21,4] = 1
$arr[21,5] = $true
$arr[22,0] = 'Sergipe'
$arr[22,1] = '01/01/2016'
$arr[22,2] = 'Homicídio doloso'
$arr[22,3] = 58.37831652311978
$arr[22,4] = 1
$arr[22,5] = $true
$arr[23,0] = 'Sergipe'
$arr[23,1] = '01/01/2017'
$arr[23,2] = 'Homicídio doloso'
$arr[23,3] = 48.77581995210135
$arr[23,4] = 6
$arr[23,5] = $false
$arr[24,0] = 'Sergipe'
$arr[24,1] = '01/01/2018'
$arr[24,2] = 'Homicídio doloso'
$arr[24,3] = 40.90755069112692
$arr[24,4] = 7
$arr[24,5] = $false
$arr[25,0] = 'Sergipe'
$arr[25,1] = '01/01/2019'
$arr[25,2] = 'Homicídio doloso'
$arr[25,3] = 32.66776922200251
$arr[25,4] = 4
$arr[25,5] = $false
$arr[26,0] = 'Sergipe'
$arr[26,1] = '01/01/2020'
$arr[26,2] = 'Homicídio doloso'
$arr[26,3] = 32.21171592285057
$arr[26,4] = 6
$arr[26,5] = $false
$arr[27,0] = 'Sergipe'
$arr[27,1] = '01/01/2021'
$arr[27,2] = 'Homicídio doloso'
$arr[27,3] = 23.17538722565815
$arr[27,4] = 14
$arr[27,5] = $false
$arr[28,0] = 'Sergipe'
$arr[28,1] = '01/01/2022'
$arr[28,2] = 'Homicídio doloso'
$arr[28,3] = 23.75055135208496
$arr[28,4] = 13
$arr[28,5] = $false
$arr[29,0] = 'Sergipe'
$arr[29,1] = '01/01/2023'
$arr[29,2] = 'Homicídio doloso'
$arr[29,3] = 18.30463713266065
$arr[29,4] = 19
$arr[29,5] = $false
$arr[30,0] = 'Sergipe'
$arr[30,1] = '01/01/2024'
$arr[30,2] = 'Homicídio doloso'
$arr[30,3] = 14.82574116177014
$arr[30,4] = 17
$arr[30,5] = $false

$ws.Range("A1:F31").Value = $arr

# Safety net: re-assert column B as literal text in case the bulk array
# write above ever got reinterpreted as a date by Excel's auto-detection.
$ws.Range("B2:B31").NumberFormat = "@"
for ($r = 1; $r -lt 31; $r++) {
    $ws.Cells.Item($r + 1, 2).Value = $arr[$r, 1]
}

# Give the new F1 header the same style as the existing header cells (bold,
# bordered, centered) by copying formatting from E1.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F1").Value = 'Faltam dados para todos os Estados'

$ws.Range("A1").Select()
